# inventor stls - make arm fit within walls
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# armh (B27): was "=B4 - 20" -> now "=B3 - 20" (arm height now derives from
# wall2_thick instead of wall3_thick so the arm fits within the walls)
$ws.Range("B27").Formula = "=B3 - 20"

# armhornthick (B34): was "=B32+2+2" -> now "=B32+2" (one less clearance step)
$ws.Range("B34").Formula = "=B32+2"

# Match the saved view state: scrolled/selected around the edited cells.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B33").Select()
